# Apply updated cryptocurrency price/volume data (and the two rank swaps)
# exactly as captured by the source diff, cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '97.390.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.347.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '652.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.38'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -10.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.416'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -11.48%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.01'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.341.10'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.208'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.11'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '97.094.06'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000253'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.967.11'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.339.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.533'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +20.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.63'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '498.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.32'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000199'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '93.54'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -9.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.538.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.146'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.994'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.190'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.48'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +16.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.546'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '28.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.66'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.37%  '
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '513.40'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.148'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '24.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.840'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.76'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.60%  '
$ws.Range('B46').Value = 'MantraDAO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.65'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.48%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0415'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.53'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.43'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.13'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.57%  '
